$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Type" column (C) for "École nationale supérieure des mines d'Albi" (row 3):
# it was incorrectly saved as "U" (Université) when another row's edit ("EI" for
# École nationale supérieure des mines d'Alès) leaked into it. The actual fix is
# simply to set C3 to "EI".
$ws.Range("C3").Value = "EI"
